$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2"); $c.NumberFormat = "@"; $c.Value = "30.387.35"; $c.Style = "Normal"
$ws.Range("E2").Value = "  +0.68%  "

# Row 3
$c = $ws.Range("D3"); $c.NumberFormat = "@"; $c.Value = "1.879.49"; $c.Style = "Normal"
$ws.Range("E3").Value = "  +0.76%  "

# Row 4
$c = $ws.Range("D4"); $c.NumberFormat = "@"; $c.Value = "1.000"; $c.Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$c = $ws.Range("D5"); $c.NumberFormat = "@"; $c.Value = "244.30"; $c.Style = "Normal"
$ws.Range("E5").Value = "  +4.12%  "

# Row 6
$c = $ws.Range("D6"); $c.NumberFormat = "@"; $c.Value = "0.9998"; $c.Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "

# Row 7
$c = $ws.Range("D7"); $c.NumberFormat = "@"; $c.Value = "0.4773"; $c.Style = "Normal"
$ws.Range("E7").Value = "  +1.81%  "

# Row 8
$c = $ws.Range("D8"); $c.NumberFormat = "@"; $c.Value = "0.2885"; $c.Style = "Normal"
$ws.Range("E8").Value = "  +1.21%  "

# Row 9
$c = $ws.Range("D9"); $c.NumberFormat = "@"; $c.Value = "0.06521"; $c.Style = "Normal"
$ws.Range("E9").Value = "  -0.47%  "

# Row 10
$c = $ws.Range("D10"); $c.NumberFormat = "@"; $c.Value = "21.41"; $c.Style = "Normal"
$ws.Range("E10").Value = "  +0.10%  "

# Row 11
$c = $ws.Range("D11"); $c.NumberFormat = "@"; $c.Value = "0.07762"; $c.Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "

# Row 12
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Range("D12"); $c.NumberFormat = "@"; $c.Value = "0.7379"; $c.Style = "Normal"
$ws.Range("E12").Value = "  +6.89%  "

# Row 13
$c = $ws.Range("D13"); $c.NumberFormat = "@"; $c.Value = "96.40"; $c.Style = "Normal"
$ws.Range("E13").Value = "  +0.63%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c = $ws.Range("D14"); $c.NumberFormat = "@"; $c.Value = "1.875.35"; $c.Style = "Normal"
$ws.Range("E14").Value = "  +0.54%  "

# Row 15
$c = $ws.Range("D15"); $c.NumberFormat = "@"; $c.Value = "5.128"; $c.Style = "Normal"
$ws.Range("E15").Value = "  +0.83%  "

# Row 16
$c = $ws.Range("D16"); $c.NumberFormat = "@"; $c.Value = "276.23"; $c.Style = "Normal"
$ws.Range("E16").Value = "  +3.74%  "

# Row 17
$c = $ws.Range("D17"); $c.NumberFormat = "@"; $c.Value = "30.375.00"; $c.Style = "Normal"
$ws.Range("E17").Value = "  +0.67%  "

# Row 18
$c = $ws.Range("D18"); $c.NumberFormat = "@"; $c.Value = "13.40"; $c.Style = "Normal"
$ws.Range("E18").Value = "  -1.92%  "

# Row 19
$c = $ws.Range("D19"); $c.NumberFormat = "@"; $c.Value = "0.000007549"; $c.Style = "Normal"
$ws.Range("E19").Value = "  -2.09%  "

# Row 20
$c = $ws.Range("D20"); $c.NumberFormat = "@"; $c.Value = "0.9996"; $c.Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "

# Row 21
$c = $ws.Range("D21"); $c.NumberFormat = "@"; $c.Value = "2.124.00"; $c.Style = "Normal"
$ws.Range("E21").Value = "  -0.17%  "

# Row 22
$c = $ws.Range("D22"); $c.NumberFormat = "@"; $c.Value = "1.0000"; $c.Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "

# Row 23
$c = $ws.Range("D23"); $c.NumberFormat = "@"; $c.Value = "5.233"; $c.Style = "Normal"
$ws.Range("E23").Value = "  -0.34%  "

# Row 24
$c = $ws.Range("D24"); $c.NumberFormat = "@"; $c.Value = "6.171"; $c.Style = "Normal"
$ws.Range("E24").Value = "  +0.32%  "

# Row 25
$c = $ws.Range("D25"); $c.NumberFormat = "@"; $c.Value = "9.260"; $c.Style = "Normal"
$ws.Range("E25").Value = "  -2.30%  "

# Row 26
$c = $ws.Range("D26"); $c.NumberFormat = "@"; $c.Value = "163.54"; $c.Style = "Normal"
$ws.Range("E26").Value = "  -1.57%  "

# Row 27
$c = $ws.Range("D27"); $c.NumberFormat = "@"; $c.Value = "18.95"; $c.Style = "Normal"
$ws.Range("E27").Value = "  +1.58%  "

# Row 28
$c = $ws.Range("D28"); $c.NumberFormat = "@"; $c.Value = "1.955"; $c.Style = "Normal"
$ws.Range("E28").Value = "  +1.15%  "

# Row 29
$c = $ws.Range("D29"); $c.NumberFormat = "@"; $c.Value = "1.370"; $c.Style = "Normal"
$ws.Range("E29").Value = "  +0.22%  "

# Row 30
$c = $ws.Range("D30"); $c.NumberFormat = "@"; $c.Value = "0.09955"; $c.Style = "Normal"
$ws.Range("E30").Value = "  +0.39%  "

# Row 31
$c = $ws.Range("D31"); $c.NumberFormat = "@"; $c.Value = "1.509"; $c.Style = "Normal"
$ws.Range("E31").Value = "  +3.41%  "

# Row 32
$c = $ws.Range("D32"); $c.NumberFormat = "@"; $c.Value = "4.313"; $c.Style = "Normal"
$ws.Range("E32").Value = "  -0.94%  "

# Row 33
$c = $ws.Range("D33"); $c.NumberFormat = "@"; $c.Value = "4.083"; $c.Style = "Normal"
$ws.Range("E33").Value = "  +1.08%  "

# Row 34
$c = $ws.Range("D34"); $c.NumberFormat = "@"; $c.Value = "0.04750"; $c.Style = "Normal"
$ws.Range("E34").Value = "  +0.39%  "

# Row 35
$c = $ws.Range("D35"); $c.NumberFormat = "@"; $c.Value = "1.123"; $c.Style = "Normal"
$ws.Range("E35").Value = "  -0.36%  "

# Row 36
$c = $ws.Range("D36"); $c.NumberFormat = "@"; $c.Value = "0.6966"; $c.Style = "Normal"
$ws.Range("E36").Value = "  -0.13%  "

# Row 37
$c = $ws.Range("D37"); $c.NumberFormat = "@"; $c.Value = "2.717"; $c.Style = "Normal"
$ws.Range("E37").Value = "  -0.04%  "

# Row 38
$c = $ws.Range("D38"); $c.NumberFormat = "@"; $c.Value = "0.01859"; $c.Style = "Normal"
$ws.Range("E38").Value = "  +0.10%  "

# Row 39
$c = $ws.Range("D39"); $c.NumberFormat = "@"; $c.Value = "2.754"; $c.Style = "Normal"
$ws.Range("E39").Value = "  -0.70%  "

# Row 40
$c = $ws.Range("D40"); $c.NumberFormat = "@"; $c.Value = "6.281"; $c.Style = "Normal"
$ws.Range("E40").Value = "  -0.35%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Range("D41"); $c.NumberFormat = "@"; $c.Value = "0.4175"; $c.Style = "Normal"
$ws.Range("E41").Value = "  +0.89%  "

# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D42"); $c.NumberFormat = "@"; $c.Value = "69.42"; $c.Style = "Normal"
$ws.Range("E42").Value = "  -3.87%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D43"); $c.NumberFormat = "@"; $c.Value = "1.908"; $c.Style = "Normal"
$ws.Range("E43").Value = "  -1.09%  "

# Row 44
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D44"); $c.NumberFormat = "@"; $c.Value = "0.8409"; $c.Style = "Normal"
$ws.Range("E44").Value = "  +0.76%  "

# Row 45
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D45"); $c.NumberFormat = "@"; $c.Value = "0.9995"; $c.Style = "Normal"
$ws.Range("E45").Value = "  -0.02%  "

# Row 46
$c = $ws.Range("D46"); $c.NumberFormat = "@"; $c.Value = "101.76"; $c.Style = "Normal"
$ws.Range("E46").Value = "  -1.16%  "

# Row 47
$c = $ws.Range("D47"); $c.NumberFormat = "@"; $c.Value = "9.265"; $c.Style = "Normal"
$ws.Range("E47").Value = "  +1.25%  "

# Row 48
$c = $ws.Range("D48"); $c.NumberFormat = "@"; $c.Value = "7.097"; $c.Style = "Normal"
$ws.Range("E48").Value = "  -0.24%  "

# Row 49
$c = $ws.Range("D49"); $c.NumberFormat = "@"; $c.Value = "35.17"; $c.Style = "Normal"
$ws.Range("E49").Value = "  +1.75%  "

# Row 50
$c = $ws.Range("D50"); $c.NumberFormat = "@"; $c.Value = "913.31"; $c.Style = "Normal"
$ws.Range("E50").Value = "  -5.51%  "

# Row 51
$c = $ws.Range("D51"); $c.NumberFormat = "@"; $c.Value = "0.05593"; $c.Style = "Normal"
$ws.Range("E51").Value = "  -0.83%  "
